# Adjust studio example params
#
# "optimize" sheet: update several parameter values.
# "studio" sheet: drop the wrap-text formatting (and the oversized row
# height it forced) on the notes cell.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# optimize sheet - parameter value updates
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("optimize")

# A scratch cell used purely to mint new text values without Excel's
# automatic number/boolean type coercion (and without picking up any
# extra cell-style/number-format baggage): build the value with a
# formula that evaluates to text, then paste-special just the value
# into the destination cell.
$scratch = $ws.Range("Z1")

function Set-TextValue($range, $text) {
    $scratch.Formula = "=""" + $text + """"
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.ClearContents()
}

Set-TextValue $ws.Range("B2") "10.000000"   # ending_regularization
Set-TextValue $ws.Range("B3") "0.050000"    # initial_relative_weighting
Set-TextValue $ws.Range("B5") "false"       # multiscale
Set-TextValue $ws.Range("B7") "10.000000"   # normals_strength
Set-TextValue $ws.Range("B9") "1000"        # optimization_iterations
Set-TextValue $ws.Range("B14") "100.000000" # starting_regularization
Set-TextValue $ws.Range("B15") "0"          # use_normals

# ---------------------------------------------------------------------
# studio sheet - remove the wrap-text style + custom row height from
# the notes value cell (text content itself is unchanged)
# ---------------------------------------------------------------------
$studio = $wb.Worksheets.Item("studio")
$studio.Range("B4").Style = "Normal"
$studio.Rows.Item(4).AutoFit()
